$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.379.74'
$ws.Range("E2").Value = '  -3.36%  '
$ws.Range("D3").Value = '1.650.29'
$ws.Range("E3").Value = '  -3.91%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("E5").Value = '  -2.22%  '
$ws.Range("E6").Value = '  -2.31%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.99'
$ws.Range("E8").Value = '  -0.85%  '
$ws.Range("E9").Value = '  -1.81%  '
$ws.Range("E10").Value = '  -2.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0875'
$ws.Range("E11").Value = '  -1.63%  '
$ws.Range("D12").Value = '1.886.44'
$ws.Range("E12").Value = '  -3.70%  '
$ws.Range("D13").Value = '1.650.39'
$ws.Range("E13").Value = '  -3.88%  '
$ws.Range("E14").Value = '  -2.68%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.563'
$ws.Range("E15").Value = '  +0.18%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.60'
$ws.Range("E16").Value = '  -2.81%  '
$ws.Range("D17").Value = '27.365.70'
$ws.Range("E17").Value = '  -3.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '234.61'
$ws.Range("E18").Value = '  -7.67%  '
$ws.Range("E19").Value = '  -2.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.51'
$ws.Range("E20").Value = '  -3.59%  '
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.41'
$ws.Range("E22").Value = '  -3.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.23'
$ws.Range("E23").Value = '  -4.06%  '
$ws.Range("E24").Value = '  -0.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.63'
$ws.Range("E25").Value = '  -1.53%  '
$ws.Range("E26").Value = '  -3.34%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.04'
$ws.Range("E27").Value = '  -3.50%  '
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("E29").Value = '  -2.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0498'
$ws.Range("E30").Value = '  -2.86%  '
$ws.Range("E31").Value = '  -0.88%  '
$ws.Range("E32").Value = '  -3.14%  '
$ws.Range("D33").Value = '1.452.33'
$ws.Range("E33").Value = '  -1.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.07'
$ws.Range("E34").Value = '  -4.22%  '
$ws.Range("E35").Value = '  -4.66%  '
$ws.Range("E36").Value = '  -0.19%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.910'
$ws.Range("E37").Value = '  -6.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.570'
$ws.Range("E38").Value = '  -4.83%  '
$ws.Range("E39").Value = '  -3.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.04'
$ws.Range("E40").Value = '  -0.58%  '
$ws.Range("E41").Value = '  +0.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '65.51'
$ws.Range("E42").Value = '  -5.86%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.44'
$ws.Range("E43").Value = '  -3.74%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.787'
$ws.Range("E45").Value = '  -2.14%  '
$ws.Range("D46").Value = '1.792.34'
$ws.Range("E46").Value = '  -3.81%  '
$ws.Range("E47").Value = '  -1.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.23'
$ws.Range("E48").Value = '  -2.07%  '
$ws.Range("D49").Value = '0.0₆0106'
$ws.Range("E49").Value = '  -6.79%  '
$ws.Range("E50").Value = '  -2.23%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.78'
$ws.Range("E51").Value = '  -3.58%  '
